$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 3203090.2
$ws.Range("J17").Value = 3203090.2
$ws.Range("L17").Value = 9609270.600000001
$ws.Range("N17").Value = -9609606.600000001

$ws.Range("H40").Value = 913595.75
$ws.Range("I40").Value = 1114839.4
$ws.Range("J40").Value = 7999.5
$ws.Range("K40").Value = 1114839.4
$ws.Range("L40").Value = 7999.5
$ws.Range("M40").Value = -1114664.4
$ws.Range("N40").Value = -8349.5

$ws.Range("H58").Value = 3081.4285
$ws.Range("I58").Value = 314
$ws.Range("K58").Value = 942
$ws.Range("M58").Value = -792

$ws.Range("H62").Value = 3213.5
$ws.Range("I62").Value = 2848.2856
$ws.Range("J62").Value = 4065.6667
$ws.Range("K62").Value = 2848.2856
$ws.Range("L62").Value = 4065.6667
$ws.Range("M62").Value = -2224.2856
$ws.Range("N62").Value = -5313.6667

$ws.Range("H65").Value = 3213.5
$ws.Range("I65").Value = 2848.2856
$ws.Range("J65").Value = 4065.6667
$ws.Range("K65").Value = 14241.428
$ws.Range("L65").Value = 20328.3335
$ws.Range("M65").Value = -11121.428
$ws.Range("N65").Value = -26568.3335

$ws.Range("H132").Value = 1774.0851
$ws.Range("I132").Value = 1764.0222
$ws.Range("K132").Value = 5292.0666
$ws.Range("M132").Value = -2762.0666

$ws.Range("H137").Value = 4423.227
$ws.Range("I137").Value = 2391.6584
$ws.Range("K137").Value = 7174.975199999999
$ws.Range("M137").Value = -4624.975199999999

$ws.Range("H138").Value = 35720044
$ws.Range("I138").Value = 1112.6666
$ws.Range("J138").Value = 62509240
$ws.Range("K138").Value = 3337.9998
$ws.Range("L138").Value = 187527720
$ws.Range("M138").Value = 1802.0002
$ws.Range("N138").Value = -187538000


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 742.13
$ws.Range("I32").Value = 722.8763
$ws.Range("J32").Value = 1364.6666
$ws.Range("K32").Value = 722.8763
$ws.Range("L32").Value = 1364.6666
$ws.Range("M32").Value = -435.8763
$ws.Range("N32").Value = -1938.6666

$ws.Range("H45").Value = 2769.5862
$ws.Range("I45").Value = 2395.2727
$ws.Range("K45").Value = 2395.2727
$ws.Range("M45").Value = -2018.2727

$ws.Range("H61").Value = 22729594
$ws.Range("I61").Value = 23811478
$ws.Range("K61").Value = 23811478
$ws.Range("M61").Value = -23811266

$ws.Range("H74").Value = 66742700
$ws.Range("I74").Value = 111236504
$ws.Range("J74").Value = 2004.1666
$ws.Range("K74").Value = 111236504
$ws.Range("L74").Value = 2004.1666
$ws.Range("M74").Value = -111235630
$ws.Range("N74").Value = -3752.1666

$ws.Range("H77").Value = 66742700
$ws.Range("I77").Value = 111236504
$ws.Range("J77").Value = 2004.1666
$ws.Range("K77").Value = 556182520
$ws.Range("L77").Value = 10020.833
$ws.Range("M77").Value = -556178152
$ws.Range("N77").Value = -18756.833

$ws.Range("H110").Value = 25787.5
$ws.Range("I110").Value = 29950.75
$ws.Range("K110").Value = 29950.75
$ws.Range("M110").Value = -27905.75

$ws.Range("H132").Value = 27849040
$ws.Range("I132").Value = 5084.355
$ws.Range("K132").Value = 15253.065
$ws.Range("M132").Value = -12723.065

$ws.Range("H136").Value = 22729594
$ws.Range("I136").Value = 23811478
$ws.Range("K136").Value = 71434434
$ws.Range("M136").Value = -71431884


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H134").Value = 1540347.5
$ws.Range("I134").Value = 1696403.4
$ws.Range("K134").Value = 5089210.199999999
$ws.Range("M134").Value = -5086675.199999999


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 38465508
$ws.Range("I31").Value = 2845.6155
$ws.Range("K31").Value = 2845.6155
$ws.Range("M31").Value = -2550.6155

$ws.Range("H34").Value = 38465508
$ws.Range("I34").Value = 2845.6155
$ws.Range("K34").Value = 2845.6155
$ws.Range("M34").Value = -2643.6155

$ws.Range("H59").Value = 43333.332
$ws.Range("I59").Value = 43333.332
$ws.Range("K59").Value = 43333.332
$ws.Range("M59").Value = -42188.332

$ws.Range("H122").Value = 1738.2
$ws.Range("I122").Value = 1163.7778
$ws.Range("J122").Value = 2599.8333
$ws.Range("K122").Value = 3491.3334
$ws.Range("L122").Value = 7799.499899999999
$ws.Range("M122").Value = -1041.3334
$ws.Range("N122").Value = -12699.4999

$ws.Range("H132").Value = 71267.42999999999
$ws.Range("I132").Value = 128345.25
$ws.Range("K132").Value = 385035.75
$ws.Range("M132").Value = -382505.75

$ws.Range("H134").Value = 2325.4707
$ws.Range("I134").Value = 2169.25
$ws.Range("K134").Value = 6507.75
$ws.Range("M134").Value = -3972.75


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H9").Value = 650
$ws.Range("I9").Value = 437.5
$ws.Range("J9").Value = 1500
$ws.Range("K9").Value = 1312.5
$ws.Range("L9").Value = 4500
$ws.Range("M9").Value = -1088.5
$ws.Range("N9").Value = -4948

$ws.Range("H49").Value = 5500
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H58").Value = 4748.75
$ws.Range("J58").Value = 5999.6665
$ws.Range("L58").Value = 17998.9995
$ws.Range("N58").Value = -18254.9995

$ws.Range("H68").Value = 2356.4285
$ws.Range("I68").Value = 1999.5
$ws.Range("K68").Value = 5998.5
$ws.Range("M68").Value = -5187.5

$ws.Range("H71").Value = 2356.4285
$ws.Range("I71").Value = 1999.5
$ws.Range("K71").Value = 17995.5
$ws.Range("M71").Value = -13939.5

$ws.Range("H122").Value = 97
$ws.Range("I122").Value = 97
$ws.Range("J122").Value = 97
$ws.Range("K122").Value = 873
$ws.Range("L122").Value = 873
$ws.Range("M122").Value = 1577
$ws.Range("N122").Value = -5773

$ws.Range("H141").Value = 6573.773
$ws.Range("I141").Value = 4699.467
$ws.Range("J141").Value = 10590.143
$ws.Range("K141").Value = 14098.401
$ws.Range("L141").Value = 31770.429
$ws.Range("M141").Value = -8918.400999999998
$ws.Range("N141").Value = -42130.429


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 2395.0454
$ws.Range("J80").Value = 2187.6667
$ws.Range("L80").Value = 2187.6667
$ws.Range("N80").Value = -4183.6667

$ws.Range("H83").Value = 2395.0454
$ws.Range("J83").Value = 2187.6667
$ws.Range("L83").Value = 10938.3335
$ws.Range("N83").Value = -20922.3335

$ws.Range("H107").Value = 605
$ws.Range("I107").Value = 532.9167
$ws.Range("K107").Value = 532.9167
$ws.Range("M107").Value = 1387.0833

$ws.Range("H132").Value = 4434.6665
$ws.Range("I132").Value = 5019.3076
$ws.Range("J132").Value = 3484.625
$ws.Range("K132").Value = 15057.9228
$ws.Range("L132").Value = 10453.875
$ws.Range("M132").Value = -12527.9228
$ws.Range("N132").Value = -15513.875


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 2136.9565
$ws.Range("I46").Value = 803.5714
$ws.Range("J46").Value = 4211.1113
$ws.Range("K46").Value = 803.5714
$ws.Range("L46").Value = 4211.1113
$ws.Range("M46").Value = -615.5714
$ws.Range("N46").Value = -4587.1113


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H14").Value = 31555.238
$ws.Range("I14").Value = 100532
$ws.Range("K14").Value = 100532
$ws.Range("M14").Value = -100364

$ws.Range("H46").Value = 126400
$ws.Range("J46").Value = 126400
$ws.Range("L46").Value = 126400
$ws.Range("N46").Value = -126862

$ws.Range("H113").Value = 526.5454999999999
$ws.Range("I113").Value = 351.44446
$ws.Range("K113").Value = 1054.33338
$ws.Range("M113").Value = 1115.66662

$ws.Range("H126").Value = 9027.786
$ws.Range("I126").Value = 9795
$ws.Range("J126").Value = 4424.5
$ws.Range("K126").Value = 29385
$ws.Range("L126").Value = 13273.5
$ws.Range("M126").Value = -26915
$ws.Range("N126").Value = -18213.5

$ws.Range("H132").Value = 2273.5454
$ws.Range("I132").Value = 2113.9355
$ws.Range("K132").Value = 6341.806500000001
$ws.Range("M132").Value = -3811.806500000001

$ws.Range("H134").Value = 126400
$ws.Range("J134").Value = 126400
$ws.Range("L134").Value = 379200
$ws.Range("N134").Value = -384270

